# Error Calculations and Plots
# This edit removes two sample rows (RM 232 and SC 92) from the missing-data
# worksheet and updates a handful of cells that represent newly
# imputed / newly-blanked "missing data" values for the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from the dataset ---
# Delete from the bottom up so row numbers for the earlier deletion stay valid.
$ws.Rows(28).Delete()   # was "SC 92"
$ws.Rows(26).Delete()   # was "RM 232"

# --- After the deletions, patch individual cells to their new (imputed /
#     newly-missing) values ---

# Row 5  (RM 14): F5 becomes missing
$ws.Range("F5").Value = ""

# Row 11 (RM 58): F11 gets an imputed value
$ws.Range("F11").Value = 17.65

# Row 19 (RM 125): E19 gets an imputed value, F19 becomes missing
$ws.Range("E19").Value = -6.5
$ws.Range("F19").Value = ""

# Row 21 (RM 135): E21 becomes missing
$ws.Range("E21").Value = ""

# Row 23 (RM 140): E23 and F23 get imputed values
$ws.Range("E23").Value = -7
$ws.Range("F23").Value = 16.48

# Row 25 (RM 145): F25 gets an imputed value
$ws.Range("F25").Value = 16.6

# Row 26 (SC 5, shifted up from row 27): B26 becomes missing
$ws.Range("B26").Value = ""

# Row 27 (SC 101, shifted up from row 29): B27 gets an imputed value,
# E27 and F27 become missing
$ws.Range("B27").Value = -20.4
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""

# Row 29 (SC 119, shifted up from row 31): B29 becomes missing, F29 becomes missing
$ws.Range("B29").Value = ""
$ws.Range("F29").Value = ""

# Row 33 (SC 232, shifted up from row 35): E33 and F33 get imputed values
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
